$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.266.26'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '1.688.32'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D5").Value = '217.26'
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").Value = '0.5347'
$ws.Range("E6").Value = '  +1.64%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.2711'
$ws.Range("E8").Value = '  +0.68%  '
$ws.Range("D9").Value = '0.06397'
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("D10").Value = '21.62'
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("D11").Value = '0.07675'
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("D12").Value = '1.682.39'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = '4.519'
$ws.Range("E13").Value = '  -0.04%  '
$ws.Range("D14").Value = '0.5765'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").Value = '0.000008316'
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("D16").Value = '66.46'
$ws.Range("E16").Value = '  +2.61%  '
$ws.Range("D17").Value = '26.294.93'
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = '4.874'
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("D20").Value = '10.83'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").Value = '190.49'
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = '6.232'
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '148.43'
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("D25").Value = '0.1283'
$ws.Range("E25").Value = '  +2.79%  '
$ws.Range("D26").Value = '7.825'
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").Value = '15.79'
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '0.06172'
$ws.Range("E28").Value = '  -5.01%  '
$ws.Range("D29").Value = '1.375'
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").Value = '1.323'
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("D31").Value = '3.588'
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").Value = '3.576'
$ws.Range("E32").Value = '  -0.50%  '
$ws.Range("D33").Value = '1.676'
$ws.Range("E33").Value = '  +0.70%  '
$ws.Range("D34").Value = '1.026'
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").Value = '2.758'
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("D38").Value = '0.01643'
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("E39").Value = '  -4.39%  '
$ws.Range("D40").Value = '1.104.95'
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").Value = '0.8781'
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").Value = '1.012'
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("D43").Value = '100.73'
$ws.Range("E43").Value = '  +0.31%  '
$ws.Range("D44").Value = '1.839.59'
$ws.Range("D45").Value = '0.00000000112'
$ws.Range("E45").Value = '  +3.55%  '
$ws.Range("D46").Value = '57.56'
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '1.002'
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '8.098'
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("D50").Value = '0.4298'
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '6.048'
$ws.Range("E51").Value = '  -0.45%  '

Write-Host "Applied all changes"
